$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 754 (pushes the existing 754:785 block down to 757:788)
$ws.Range("A754:T756").EntireRow.Insert()

# --- Row 754 ---
$ws.Cells.Item(754, 1).Value = 8
$ws.Cells.Item(754, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(754, 3).Value = "Coquimbo"
$ws.Cells.Item(754, 4).Value = 44747
$ws.Cells.Item(754, 5).Value = 4
$ws.Cells.Item(754, 6).Value = "Fruta"
$ws.Cells.Item(754, 7).Value = 100101
$ws.Cells.Item(754, 8).Value = "Berries"
$ws.Cells.Item(754, 9).Value = 100112025
$ws.Cells.Item(754, 10).Value = "Frutilla"
$ws.Cells.Item(754, 11).Value = "Sin especificar"
$ws.Cells.Item(754, 12).Value = "Especial"
$ws.Cells.Item(754, 13).Value = 360
$ws.Cells.Item(754, 14).Value = 25000
$ws.Cells.Item(754, 15).Value = 26000
$ws.Cells.Item(754, 16).Value = 25500
$ws.Cells.Item(754, 17).Value = '$/bandeja 7 kilos'
$ws.Cells.Item(754, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(754, 19).Value = 3643
$ws.Cells.Item(754, 20).Value = 7

# --- Row 755 ---
$ws.Cells.Item(755, 1).Value = 8
$ws.Cells.Item(755, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(755, 3).Value = "Coquimbo"
$ws.Cells.Item(755, 4).Value = 44747
$ws.Cells.Item(755, 5).Value = 4
$ws.Cells.Item(755, 6).Value = "Fruta"
$ws.Cells.Item(755, 7).Value = 100101
$ws.Cells.Item(755, 8).Value = "Berries"
$ws.Cells.Item(755, 9).Value = 100112025
$ws.Cells.Item(755, 10).Value = "Frutilla"
$ws.Cells.Item(755, 11).Value = "Sin especificar"
$ws.Cells.Item(755, 12).Value = "Primera"
$ws.Cells.Item(755, 13).Value = 300
$ws.Cells.Item(755, 14).Value = 20000
$ws.Cells.Item(755, 15).Value = 21000
$ws.Cells.Item(755, 16).Value = 20500
$ws.Cells.Item(755, 17).Value = '$/bandeja 7 kilos'
$ws.Cells.Item(755, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(755, 19).Value = 2929
$ws.Cells.Item(755, 20).Value = 7

# --- Row 756 ---
$ws.Cells.Item(756, 1).Value = 8
$ws.Cells.Item(756, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(756, 3).Value = "Coquimbo"
$ws.Cells.Item(756, 4).Value = 44747
$ws.Cells.Item(756, 5).Value = 4
$ws.Cells.Item(756, 6).Value = "Fruta"
$ws.Cells.Item(756, 7).Value = 100101
$ws.Cells.Item(756, 8).Value = "Berries"
$ws.Cells.Item(756, 9).Value = 100112025
$ws.Cells.Item(756, 10).Value = "Frutilla"
$ws.Cells.Item(756, 11).Value = "Sin especificar"
$ws.Cells.Item(756, 12).Value = "Segunda"
$ws.Cells.Item(756, 13).Value = 240
$ws.Cells.Item(756, 14).Value = 16000
$ws.Cells.Item(756, 15).Value = 17000
$ws.Cells.Item(756, 16).Value = 16500
$ws.Cells.Item(756, 17).Value = '$/bandeja 7 kilos'
$ws.Cells.Item(756, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(756, 19).Value = 2357
$ws.Cells.Item(756, 20).Value = 7
